# Reattach the three "Max Pooling Layer" rounded-rectangle boxes to their
# outgoing straight arrow connectors (glue the connector's start point to
# connection site idx=3 — the right-middle handle — of the corresponding
# rounded rectangle), matching the geometry PowerPoint computes once the
# connector is glued/dragged onto that connection site.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- Connector "Straight Arrow Connector 9" (id=10) -> Rounded Rectangle 8 (id=9), site idx 3
$conn1 = $s.Shapes.Item(5)
$rect9 = $s.Shapes.Item(4)
$conn1.ConnectorFormat.BeginConnect($rect9, 3)
$conn1.Left   = 476.03875434454034
$conn1.Top    = 260.3883542168703
$conn1.Width  = 72.11653543307087
$conn1.Height = 0.46606297823622467

# ---- Connector "Straight Arrow Connector 15" (id=16) -> Rounded Rectangle 14 (id=15), site idx 3
$conn2 = $s.Shapes.Item(10)
$rect15 = $s.Shapes.Item(9)
$conn2.ConnectorFormat.BeginConnect($rect15, 3)
$conn2.Left   = 476.8542661798681
$conn2.Top    = 430.13606299212597
$conn2.Width  = 80.50496302915604
$conn2.Height = 0.23291338582677165
$conn2.Flip(1)   # msoFlipVertical -> sets flipV="1"

# ---- Connector "Straight Arrow Connector 19" (id=20) -> Rounded Rectangle 18 (id=19), site idx 3
$conn3 = $s.Shapes.Item(14)
$rect19 = $s.Shapes.Item(13)
$conn3.ConnectorFormat.BeginConnect($rect19, 3)
$conn3.Left   = 483.84456692913386
$conn3.Top    = 88.54369814859366
$conn3.Width  = 72.11653543307087
$conn3.Height = 0.46606297823622467
